# Update "want to go" counts (column F) across the four worksheets to
# reflect newly scraped totals, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(4, 6).Value  = 75
$ws1.Cells.Item(5, 6).Value  = 9354
$ws1.Cells.Item(6, 6).Value  = 9354
$ws1.Cells.Item(7, 6).Value  = 582
$ws1.Cells.Item(10, 6).Value = 261
$ws1.Cells.Item(11, 6).Value = 401
$ws1.Cells.Item(14, 6).Value = 430
$ws1.Cells.Item(15, 6).Value = 11967
$ws1.Cells.Item(16, 6).Value = 11967
$ws1.Cells.Item(24, 6).Value = 232
$ws1.Cells.Item(26, 6).Value = 24
$ws1.Cells.Item(27, 6).Value = 171
$ws1.Cells.Item(32, 6).Value = 2096
$ws1.Cells.Item(34, 6).Value = 13
$ws1.Cells.Item(35, 6).Value = 53
$ws1.Cells.Item(36, 6).Value = 2144
$ws1.Cells.Item(37, 6).Value = 989
$ws1.Cells.Item(38, 6).Value = 4188
$ws1.Cells.Item(39, 6).Value = 3619
$ws1.Cells.Item(40, 6).Value = 495
$ws1.Cells.Item(42, 6).Value = 3053
$ws1.Cells.Item(43, 6).Value = 1312
$ws1.Cells.Item(47, 6).Value = 498
$ws1.Cells.Item(48, 6).Value = 63
$ws1.Cells.Item(50, 6).Value = 123

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(9, 6).Value  = 41
$ws2.Cells.Item(10, 6).Value = 16
$ws2.Cells.Item(20, 6).Value = 7
$ws2.Cells.Item(24, 6).Value = 35

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 6).Value = 49

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(9, 6).Value  = 9354
$ws4.Cells.Item(10, 6).Value = 582
$ws4.Cells.Item(13, 6).Value = 261
$ws4.Cells.Item(14, 6).Value = 401
$ws4.Cells.Item(17, 6).Value = 11967
$ws4.Cells.Item(18, 6).Value = 11967
$ws4.Cells.Item(22, 6).Value = 16
$ws4.Cells.Item(23, 6).Value = 49
$ws4.Cells.Item(28, 6).Value = 171
$ws4.Cells.Item(33, 6).Value = 2096
$ws4.Cells.Item(35, 6).Value = 13
$ws4.Cells.Item(36, 6).Value = 53
$ws4.Cells.Item(38, 6).Value = 2144
$ws4.Cells.Item(39, 6).Value = 989
$ws4.Cells.Item(41, 6).Value = 7
$ws4.Cells.Item(42, 6).Value = 3619
$ws4.Cells.Item(43, 6).Value = 3053
$ws4.Cells.Item(45, 6).Value = 1312
$ws4.Cells.Item(48, 6).Value = 35
$ws4.Cells.Item(49, 6).Value = 498
$ws4.Cells.Item(50, 6).Value = 63
